$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on both the "展览" sheet
# and the "全部类型" sheet (duplicate data), matching the site regeneration.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1691
    $ws.Range("F6").Value = 466
    $ws.Range("F9").Value = 611
}
